# Scheduled runner update: refresh market-price driven profit calculations
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# per-crafting-class "Profits" sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1429672
$ws.Range("J43").Value = 2500925.5
$ws.Range("L43").Value = 2500925.5
$ws.Range("N43").Value = -2501063.5
$ws.Range("H86").Value = 3212.5
$ws.Range("I86").Value = 2876.6
$ws.Range("K86").Value = 2876.6
$ws.Range("M86").Value = -1753.6
$ws.Range("H89").Value = 3212.5
$ws.Range("I89").Value = 2876.6
$ws.Range("K89").Value = 14383
$ws.Range("M89").Value = -8767
$ws.Range("H116").Value = 2783.2173
$ws.Range("I116").Value = 3135.7144
$ws.Range("J116").Value = 2234.889
$ws.Range("K116").Value = 3135.7144
$ws.Range("L116").Value = 2234.889
$ws.Range("M116").Value = 306.2856000000002
$ws.Range("N116").Value = -9118.888999999999
$ws.Range("H132").Value = 3815.4167
$ws.Range("I132").Value = 4025.9092
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 12077.7276
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -9547.7276
$ws.Range("N132").Value = -9560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 40824.117
$ws.Range("I2").Value = 54504.58
$ws.Range("J2").Value = 3691.4285
$ws.Range("K2").Value = 54504.58
$ws.Range("L2").Value = 3691.4285
$ws.Range("M2").Value = -54391.58
$ws.Range("N2").Value = -3917.4285
$ws.Range("H110").Value = 488.91666
$ws.Range("I110").Value = 407.44446
$ws.Range("J110").Value = 733.3333
$ws.Range("K110").Value = 407.44446
$ws.Range("L110").Value = 733.3333
$ws.Range("M110").Value = 1637.55554
$ws.Range("N110").Value = -4823.3333
$ws.Range("H116").Value = 40824.117
$ws.Range("I116").Value = 54504.58
$ws.Range("J116").Value = 3691.4285
$ws.Range("K116").Value = 54504.58
$ws.Range("L116").Value = 3691.4285
$ws.Range("M116").Value = -52210.58
$ws.Range("N116").Value = -8279.4285
$ws.Range("H132").Value = 2110.3103
$ws.Range("I132").Value = 1238.1904
$ws.Range("J132").Value = 4399.625
$ws.Range("K132").Value = 3714.5712
$ws.Range("L132").Value = 13198.875
$ws.Range("M132").Value = -1184.5712
$ws.Range("N132").Value = -18258.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 40824.117
$ws.Range("I3").Value = 54504.58
$ws.Range("J3").Value = 3691.4285
$ws.Range("K3").Value = 54504.58
$ws.Range("L3").Value = 3691.4285
$ws.Range("M3").Value = -54390.58
$ws.Range("N3").Value = -3919.4285
$ws.Range("H105").Value = 5002.857
$ws.Range("I105").Value = 6505
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 6505
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -4758
$ws.Range("N105").Value = -6494
$ws.Range("H134").Value = 36691.586
$ws.Range("I134").Value = 2142.5454
$ws.Range("J134").Value = 145274.28
$ws.Range("K134").Value = 6427.6362
$ws.Range("L134").Value = 435822.84
$ws.Range("M134").Value = -3892.6362
$ws.Range("N134").Value = -440892.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2615.147
$ws.Range("I86").Value = 2780.0417
$ws.Range("J86").Value = 2219.4
$ws.Range("K86").Value = 2780.0417
$ws.Range("L86").Value = 2219.4
$ws.Range("M86").Value = -1657.0417
$ws.Range("N86").Value = -4465.4
$ws.Range("H89").Value = 2615.147
$ws.Range("I89").Value = 2780.0417
$ws.Range("J89").Value = 2219.4
$ws.Range("K89").Value = 13900.2085
$ws.Range("L89").Value = 11097
$ws.Range("M89").Value = -8284.208500000001
$ws.Range("N89").Value = -22329
$ws.Range("H109").Value = 18000
$ws.Range("J109").Value = 18000
$ws.Range("L109").Value = 18000
$ws.Range("N109").Value = -20080
$ws.Range("H132").Value = 2203.087
$ws.Range("I132").Value = 1588.2
$ws.Range("J132").Value = 3356
$ws.Range("K132").Value = 4764.6
$ws.Range("L132").Value = 10068
$ws.Range("M132").Value = -2234.6
$ws.Range("N132").Value = -15128

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 898
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 972.5
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 2917.5
$ws.Range("M92").Value = -552
$ws.Range("N92").Value = -5413.5
$ws.Range("H114").Value = 1474.826
$ws.Range("I114").Value = 1019.6667
$ws.Range("J114").Value = 1635.4706
$ws.Range("K114").Value = 3059.0001
$ws.Range("L114").Value = 4906.4118
$ws.Range("M114").Value = 194.9998999999998
$ws.Range("N114").Value = -11414.4118
$ws.Range("H131").Value = 27260.975
$ws.Range("J131").Value = 2513.3333
$ws.Range("L131").Value = 7539.999899999999
$ws.Range("N131").Value = -17619.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1492
$ws.Range("I97").Value = 1117.5
$ws.Range("J97").Value = 2490.6667
$ws.Range("K97").Value = 1117.5
$ws.Range("L97").Value = 2490.6667
$ws.Range("M97").Value = -621.5
$ws.Range("N97").Value = -3482.6667
$ws.Range("H122").Value = 1445644.5
$ws.Range("I122").Value = 1625725.1
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 4877175.300000001
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4874725.300000001
$ws.Range("N122").Value = -19900
$ws.Range("H130").Value = 50780
$ws.Range("J130").Value = 50780
$ws.Range("L130").Value = 50780
$ws.Range("N130").Value = -60820
$ws.Range("H132").Value = 2444.4324
$ws.Range("I132").Value = 1841.6428
$ws.Range("J132").Value = 4319.778
$ws.Range("K132").Value = 5524.928400000001
$ws.Range("L132").Value = 12959.334
$ws.Range("M132").Value = -2994.928400000001
$ws.Range("N132").Value = -18019.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 697.96295
$ws.Range("I61").Value = 592.6316
$ws.Range("J61").Value = 948.125
$ws.Range("K61").Value = 592.6316
$ws.Range("L61").Value = 948.125
$ws.Range("M61").Value = -390.6316
$ws.Range("N61").Value = -1352.125
$ws.Range("H93").Value = 2002
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 2004
$ws.Range("K93").Value = 2000
$ws.Range("L93").Value = 2004
$ws.Range("M93").Value = -752
$ws.Range("N93").Value = -4500
$ws.Range("H113").Value = 697.96295
$ws.Range("I113").Value = 592.6316
$ws.Range("J113").Value = 948.125
$ws.Range("K113").Value = 592.6316
$ws.Range("L113").Value = 948.125
$ws.Range("M113").Value = 1577.3684
$ws.Range("N113").Value = -5288.125
$ws.Range("H118").Value = 75000
$ws.Range("J118").Value = 75000
$ws.Range("L118").Value = 75000
$ws.Range("N118").Value = -78314
$ws.Range("H125").Value = 33843.08
$ws.Range("J125").Value = 33843.08
$ws.Range("L125").Value = 33843.08
$ws.Range("N125").Value = -43683.08
$ws.Range("H127").Value = 48320
$ws.Range("J127").Value = 48320
$ws.Range("L127").Value = 48320
$ws.Range("N127").Value = -58240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 9394.833000000001
$ws.Range("I52").Value = 2133.3333
$ws.Range("J52").Value = 16656.334
$ws.Range("K52").Value = 2133.3333
$ws.Range("L52").Value = 16656.334
$ws.Range("M52").Value = -1907.3333
$ws.Range("N52").Value = -17108.334
$ws.Range("H62").Value = 11800
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 11800
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 11800
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -13048
$ws.Range("H65").Value = 11800
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 11800
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 59000
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -65240
$ws.Range("H122").Value = 1197.2
$ws.Range("I122").Value = 995
$ws.Range("J122").Value = 1283.8572
$ws.Range("K122").Value = 2985
$ws.Range("L122").Value = 3851.5716
$ws.Range("M122").Value = -535
$ws.Range("N122").Value = -8751.571599999999
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""
